# Updated cryptos list values (Price + Volume(1h) columns), matching a GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.905.81'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '3.949.53'
$ws.Range("E3").Value = '  +2.47%  '
$ws.Range("D4").Value = "'" + '1.00'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'" + '609.12'
$ws.Range("E5").Value = '  +0.94%  '
$ws.Range("D6").Value = "'" + '169.98'
$ws.Range("E6").Value = '  +3.01%  '
$ws.Range("D7").Value = '3.949.90'
$ws.Range("E7").Value = '  +2.55%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = "'" + '0.536'
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("D11").Value = "'" + '6.48'
$ws.Range("E11").Value = '  +2.29%  '
$ws.Range("D12").Value = "'" + '0.468'
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("D13").Value = "'" + '0.0000258'
$ws.Range("E13").Value = '  +5.01%  '
$ws.Range("D14").Value = "'" + '38.06'
$ws.Range("E14").Value = '  +2.14%  '
$ws.Range("D15").Value = '4.609.59'
$ws.Range("E15").Value = '  +2.43%  '
$ws.Range("D16").Value = '3.947.29'
$ws.Range("E16").Value = '  +2.29%  '
$ws.Range("D17").Value = '69.799.91'
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("D18").Value = "'" + '7.56'
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("D19").Value = "'" + '17.45'
$ws.Range("E19").Value = '  +1.21%  '
$ws.Range("D20").Value = "'" + '0.112'
$ws.Range("E20").Value = '  -1.88%  '
$ws.Range("D21").Value = "'" + '11.10'
$ws.Range("E21").Value = '  -5.19%  '
$ws.Range("D22").Value = "'" + '500.53'
$ws.Range("E22").Value = '  +2.08%  '
$ws.Range("D23").Value = "'" + '0.738'
$ws.Range("E23").Value = '  +1.85%  '
$ws.Range("D24").Value = "'" + '0.0000168'
$ws.Range("E24").Value = '  +5.98%  '
$ws.Range("D25").Value = "'" + '85.48'
$ws.Range("E25").Value = '  +0.94%  '
$ws.Range("D26").Value = "'" + '2.30'
$ws.Range("E26").Value = '  +0.98%  '
$ws.Range("D27").Value = "'" + '12.30'
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("D28").Value = "'" + '10.29'
$ws.Range("E28").Value = '  +2.23%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").Value = "'" + '3.00'
$ws.Range("E30").Value = '  +0.42%  '
$ws.Range("D31").Value = '4.102.26'
$ws.Range("E31").Value = '  +2.38%  '
$ws.Range("D32").Value = "'" + '2.42'
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").Value = "'" + '7.87'
$ws.Range("E33").Value = '  -1.70%  '
$ws.Range("D34").Value = "'" + '32.27'
$ws.Range("E34").Value = '  -0.44%  '
$ws.Range("D35").Value = '3.921.54'
$ws.Range("E35").Value = '  +3.28%  '
$ws.Range("D36").Value = "'" + '0.108'
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").Value = "'" + '6.11'
$ws.Range("E37").Value = '  +3.16%  '
$ws.Range("B38").Value = 'Mantle'
$ws.Range("C38").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D38").Value = "'" + '1.04'
$ws.Range("E38").Value = '  +0.55%  '
$ws.Range("D39").Value = "'" + '0.141'
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("D40").Value = "'" + '3.29'
$ws.Range("E40").Value = '  +8.80%  '
$ws.Range("D41").Value = "'" + '0.999'
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").Value = "'" + '0.325'
$ws.Range("E42").Value = '  +1.30%  '
$ws.Range("D43").Value = "'" + '2.06'
$ws.Range("E43").Value = '  +2.92%  '
$ws.Range("D44").Value = "'" + '436.18'
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("D45").Value = "'" + '48.31'
$ws.Range("E45").Value = '  -0.53%  '
$ws.Range("E46").Value = '  +2.23%  '
$ws.Range("D48").Value = "'" + '0.000278'
$ws.Range("E48").Value = '  +22.48%  '
$ws.Range("D49").Value = "'" + '0.0365'
$ws.Range("E49").Value = '  +2.15%  '
$ws.Range("D50").Value = "'" + '143.11'
$ws.Range("E50").Value = '  -0.12%  '
$ws.Range("D51").Value = '2.807.56'
$ws.Range("E51").Value = '  -1.24%  '
